$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 5-7 duplicate the existing data rows (2-4) appended again.
# Columns G:K hold numeric-looking values that must stay stored as text
# (matching the rest of the sheet, which keeps every value as a string),
# so force a text number format on those columns before writing them.
$ws.Range("G5:K7").NumberFormat = "@"

# Row 5
$ws.Range("A5").Value = " Sharjah"
$ws.Range("B5").Value = " October 26 2020"
$ws.Range("C5").Value = "Kings XI won by 8 wickets (with 7 balls remaining)"
$ws.Range("D5").Value = "Kolkata Knight Riders"
$ws.Range("E5").Value = "Kings XI Punjab"
$ws.Range("F5").Value = "Varun Chakravarthy "
$ws.Range("G5").Value = "2"
$ws.Range("H5").Value = "4"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "50.00"

# Row 6
$ws.Range("A6").Value = " Sharjah"
$ws.Range("B6").Value = " October 12 2020"
$ws.Range("C6").Value = "RCB won by 82 runs"
$ws.Range("D6").Value = "Kolkata Knight Riders"
$ws.Range("E6").Value = "Royal Challengers Bangalore"
$ws.Range("F6").Value = "Varun Chakravarthy "
$ws.Range("G6").Value = "7"
$ws.Range("H6").Value = "10"
$ws.Range("I6").Value = "0"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "70.00"

# Row 7
$ws.Range("A7").Value = " Abu Dhabi"
$ws.Range("B7").Value = " October 07 2020"
$ws.Range("C7").Value = "KKR won by 10 runs"
$ws.Range("D7").Value = "Kolkata Knight Riders"
$ws.Range("E7").Value = "Chennai Super Kings"
$ws.Range("F7").Value = "Varun Chakravarthy "
$ws.Range("G7").Value = "1"
$ws.Range("H7").Value = "1"
$ws.Range("I7").Value = "0"
$ws.Range("J7").Value = "0"
$ws.Range("K7").Value = "100.00"
